# Edit slide 2 ("Icons") of the presentation:
#  - Reposition the "Picture 2" shape (small vertical nudge up).
#  - Move the "Grupo 10" group and the "Estrela de 5 pontas 11" star from the
#    front of the icon cluster to the end (on top, z-order wise), resizing /
#    repositioning them in the process.
#  - Add a new transparent "Retângulo 61" rectangle shape (a profile-picture
#    placeholder/crop frame) sitting behind the star/group, at the very end
#    of the shape stack (so it renders first among the three new additions).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$EMU = 12700.0  # points -> EMU

# ---------------------------------------------------------------------
# 1) Nudge "Picture 2" up a little.
# ---------------------------------------------------------------------
$pic2 = $s.Shapes.Item("Picture 2")
$pic2.Top = 147306 / $EMU

# ---------------------------------------------------------------------
# 2) Grab handles to the shapes that will be rebuilt at the end of the
#    z-order: the "Grupo 10" group and the "Estrela de 5 pontas 11" star.
# ---------------------------------------------------------------------
$grupo10 = $s.Shapes.Item("Grupo 10")
$estrela = $s.Shapes.Item("Estrela de 5 pontas 11")

# ---------------------------------------------------------------------
# 3) Create the new "Retângulo 61" shape by duplicating the star (so it
#    inherits the same theme shape style), then turning it into a plain
#    transparent rectangle and placing it first among the three shapes
#    that will end up on top of the stack.
# ---------------------------------------------------------------------
$novoDup = $estrela.Duplicate()
$retangulo61 = $novoDup.Item(1)
$retangulo61.Name = "Retângulo 61"
$retangulo61.AutoShapeType = 1          # msoShapeRectangle
$retangulo61.Fill.Visible = 0           # <a:noFill/>
$retangulo61.Line.Visible = 0           # <a:ln><a:noFill/></a:ln>
$retangulo61.TextFrame.TextRange.Text = ""
$retangulo61.Left = 2854683 / $EMU
$retangulo61.Top = 257459 / $EMU
$retangulo61.Width = 1447616 / $EMU
$retangulo61.Height = 1447616 / $EMU
$retangulo61.ZOrder(0)                  # msoBringToFront

# ---------------------------------------------------------------------
# 4) Resize / reposition the star, then bring it to the front (it lands
#    right after the new rectangle).
# ---------------------------------------------------------------------
$estrela.Left = -324544 / $EMU
$estrela.Top = 1132395 / $EMU
$estrela.Width = 1145359 / $EMU
$estrela.Height = 1145359 / $EMU
$estrela.ZOrder(0)                      # msoBringToFront

# ---------------------------------------------------------------------
# 5) Resize / reposition the "Grupo 10" group, then bring it to the
#    front so it becomes the very last shape in the tree.
# ---------------------------------------------------------------------
$grupo10.Left = 1045177 / $EMU
$grupo10.Top = 257459 / $EMU
$grupo10.Width = 1001883 / $EMU
$grupo10.Height = 1151589 / $EMU
$grupo10.ZOrder(0)                      # msoBringToFront
